# Update cryptos list — apply price/volume changes and re-rank the
# Aave / Cronos / TrustWalletToken block (rows 49-51).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.143.06"
$ws.Range("E2").Value = "  -1.59%  "
$ws.Range("D3").Value = "2.275.48"
$ws.Range("E3").Value = "  -1.74%  "
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").Value = "112.41"
$ws.Range("E5").Value = "  -1.93%  "
$ws.Range("D6").Value = "264.98"
$ws.Range("E6").Value = "  -2.08%  "
$ws.Range("E7").Value = "  -1.27%  "
$ws.Range("E8").Value = "  +0.25%  "
$ws.Range("D9").Value = "0.608"
$ws.Range("E9").Value = "  -2.51%  "
$ws.Range("D10").Value = "47.84"
$ws.Range("E10").Value = "  +0.21%  "
$ws.Range("E11").Value = "  -1.54%  "
$ws.Range("D12").Value = "8.78"
$ws.Range("E12").Value = "  -1.41%  "
$ws.Range("E13").Value = "  +1.35%  "
$ws.Range("D14").Value = "15.48"
$ws.Range("E14").Value = "  -1.32%  "
$ws.Range("D15").Value = "2.615.87"
$ws.Range("D16").Value = "0.855"
$ws.Range("E16").Value = "  -0.81%  "
$ws.Range("D17").Value = "2.275.09"
$ws.Range("E17").Value = "  -1.37%  "
$ws.Range("D18").Value = "43.181.78"
$ws.Range("E19").Value = "  -2.56%  "
$ws.Range("E20").Value = "  +0.88%  "
$ws.Range("D21").Value = "71.26"
$ws.Range("E21").Value = "  -2.16%  "
$ws.Range("E22").Value = "  +0.58%  "
$ws.Range("D23").Value = "231.79"
$ws.Range("E23").Value = "  -1.14%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.70"
$ws.Range("E24").Value = "  +1.66%  "
$ws.Range("D25").Value = "2.89"
$ws.Range("E25").Value = "  -1.78%  "
$ws.Range("E26").Value = "  +0.33%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.30"
$ws.Range("E27").Value = "  -1.62%  "
$ws.Range("D28").Value = "3.91"
$ws.Range("E28").Value = "  -1.35%  "
$ws.Range("D29").Value = "40.29"
$ws.Range("E29").Value = "  -5.44%  "
$ws.Range("E30").Value = "  -3.40%  "
$ws.Range("E31").Value = "  -1.68%  "
$ws.Range("D32").Value = "172.64"
$ws.Range("E32").Value = "  -2.93%  "
$ws.Range("D33").Value = "21.31"
$ws.Range("E33").Value = "  -3.08%  "
$ws.Range("D34").Value = "0.0907"
$ws.Range("E34").Value = "  -2.81%  "
$ws.Range("E35").Value = "  +2.57%  "
$ws.Range("E36").Value = "  +0.34%  "
$ws.Range("D37").Value = "4.64"
$ws.Range("E37").Value = "  -1.91%  "
$ws.Range("E38").Value = "  -1.73%  "
$ws.Range("D39").Value = "3.82"
$ws.Range("E39").Value = "  -3.90%  "
$ws.Range("D40").Value = "0.104"
$ws.Range("E40").Value = "  -6.68%  "
$ws.Range("D41").Value = "2.62"
$ws.Range("E41").Value = "  +9.10%  "
$ws.Range("D42").Value = "77.07"
$ws.Range("E42").Value = "  +8.88%  "
$ws.Range("D43").Value = "13.82"
$ws.Range("E43").Value = "  +8.00%  "
$ws.Range("E44").Value = "  -3.72%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "6.10"
$ws.Range("E45").Value = "  +2.20%  "
$ws.Range("E46").Value = "  +0.13%  "
$ws.Range("E47").Value = "  -2.01%  "
$ws.Range("E48").Value = "  -2.64%  "
$ws.Range("B49").Value = "Aave"
$ws.Range("C49").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "101.50"
$ws.Range("E49").Value = "  +0.93%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "0.0992"
$ws.Range("E50").Value = "  -1.68%  "
$ws.Range("B51").Value = "TrustWalletToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D51").Value = "1.25"
$ws.Range("E51").Value = "  +0.88%  "